$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Populate the new Neo4j query text into A2 (adds a 6th shared string,
# keeps the existing wrap-text style already applied to A2)
$ws.Range("A2").Value = 'MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE c.disease IN [''Ovarian epithelial cancer''] RETURN coalesce(c.case_id,'''') AS `Case ID` , coalesce(t.clinical_trial_designation ,'''')as `Trial Code` , coalesce(a.arm_id,'''') As `Arm` , coalesce(a.arm_drug,'''') As `Arm Treatment` , coalesce(c.disease,'''') As Diagnosis , coalesce(c.gender,'''') As Gender , coalesce(c.race,'''') As Race , coalesce(c.ethnicity,'''') As Ethnicity'

# The row grows to fit the multi-line query text
$ws.Rows.Item(2).RowHeight = 87

# Move the selection / active cell to A2 (also drops the stale topLeftCell scroll anchor)
$ws.Range("A2").Select() | Out-Null
